# Update countries & provincias Spain
# - Uzbekistan's case counts grew, moving it up two rows in the ranking
#   (it now sits where Marruecos was, shifting Marruecos and Moldavia down).
# - Refresh several countries' stats (India, Pakistan, Hong Kong, Tailandia,
#   Mauricio, Butan).
# - Update the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 07:19"

# --- Row 6: India ---
$ws.Range("B6").Value = 1805838
$ws.Range("C6").Value = 1136
$ws.Range("D6").Value = 1188389
$ws.Range("E6").Value = 579273
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 15
$ws.Range("H6").Value = 38176

# --- Row 16: Pakistan ---
$ws.Range("B16").Value = 280029
$ws.Range("C16").Value = 331
$ws.Range("D16").Value = 248873
$ws.Range("E16").Value = 25172
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 5984

# --- Rows 63-65: Uzbekistan moves up (ahead of Marruecos & Moldavia) ---
# Row 63 becomes Uzbekistan with fresh data
$ws.Range("A63").Value = "Uzbekistan"
$ws.Range("B63").Value = 25553
$ws.Range("C63").Value = 217
$ws.Range("D63").Value = 16507
$ws.Range("E63").Value = 8893
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 153

# Row 64 becomes Marruecos (previously row 63's data)
$ws.Range("A64").Value = "Marruecos"
$ws.Range("B64").Value = 25537
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 18435
$ws.Range("E64").Value = 6720
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 382

# Row 65 becomes Moldavia (previously row 64's data)
$ws.Range("A65").Value = "Moldavia"
$ws.Range("B65").Value = 25362
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 17816
$ws.Range("E65").Value = 6755
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 791

# --- Row 112: Hong Kong ---
$ws.Range("E112").Value = 1516
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 3
$ws.Range("H112").Value = 37

# --- Row 113: Tailandia ---
$ws.Range("B113").Value = 3320
$ws.Range("C113").Value = 3
$ws.Range("D113").Value = 3142
$ws.Range("E113").Value = 120

# --- Row 170: Mauricio ---
$ws.Range("D170").Value = 334
$ws.Range("E170").Value = 0

# --- Row 190: Butan ---
$ws.Range("B190").Value = 103
$ws.Range("C190").Value = 1
$ws.Range("D190").Value = 89
$ws.Range("E190").Value = 14
